$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 633.8333
$ws.Range("I38").Value = 160.6
$ws.Range("K38").Value = 481.8
$ws.Range("M38").Value = -109.8

$ws.Range("H49").Value = 575
$ws.Range("I49").Value = 300
$ws.Range("J49").Value = 850
$ws.Range("K49").Value = 900
$ws.Range("L49").Value = 2550
$ws.Range("M49").Value = -764
$ws.Range("N49").Value = -2822

$ws.Range("H64").Value = 6700
$ws.Range("I64").Value = 6500
$ws.Range("J64").Value = 6900
$ws.Range("K64").Value = 6500
$ws.Range("L64").Value = 6900
$ws.Range("M64").Value = -6252
$ws.Range("N64").Value = -7396

$ws.Range("H67").Value = 6700
$ws.Range("I67").Value = 6500
$ws.Range("J67").Value = 6900
$ws.Range("K67").Value = 6500
$ws.Range("L67").Value = 6900
$ws.Range("M67").Value = -5642
$ws.Range("N67").Value = -8616

$ws.Range("H70").Value = 754.1786
$ws.Range("I70").Value = 622.4783
$ws.Range("J70").Value = 1360
$ws.Range("K70").Value = 1867.4349
$ws.Range("L70").Value = 4080
$ws.Range("M70").Value = -1597.4349
$ws.Range("N70").Value = -4620

$ws.Range("H73").Value = 754.1786
$ws.Range("I73").Value = 622.4783
$ws.Range("J73").Value = 1360
$ws.Range("K73").Value = 1867.4349
$ws.Range("L73").Value = 4080
$ws.Range("M73").Value = -931.4349
$ws.Range("N73").Value = -5952

$ws.Range("H106").Value = 2983.25
$ws.Range("I106").Value = 3122.8572
$ws.Range("K106").Value = 3122.8572
$ws.Range("M106").Value = -2491.8572

$ws.Range("H132").Value = 1024903.5
$ws.Range("I132").Value = 3983.8206
$ws.Range("J132").Value = 5448889
$ws.Range("K132").Value = 11951.4618
$ws.Range("L132").Value = 16346667
$ws.Range("M132").Value = -9421.461800000001
$ws.Range("N132").Value = -16351727

$ws.Range("H133").Value = 24963.334
$ws.Range("J133").Value = 24963.334
$ws.Range("L133").Value = 24963.334
$ws.Range("N133").Value = -35083.334

$ws.Range("H135").Value = 41931.76
$ws.Range("I135").Value = 43512.25
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 391610.25
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -389075.25
$ws.Range("N135").Value = -41070

$ws.Range("H137").Value = 1964072
$ws.Range("I137").Value = 3335347.2
$ws.Range("J137").Value = 5107.2856
$ws.Range("K137").Value = 10006041.6
$ws.Range("L137").Value = 15321.8568
$ws.Range("M137").Value = -10003491.6
$ws.Range("N137").Value = -20421.8568

$ws.Range("H138").Value = 5130671
$ws.Range("I138").Value = 1801.9524
$ws.Range("J138").Value = 11114352
$ws.Range("K138").Value = 5405.857199999999
$ws.Range("L138").Value = 33343056
$ws.Range("M138").Value = -265.8571999999995
$ws.Range("N138").Value = -33353336

$ws.Range("H141").Value = 1272.6471
$ws.Range("I141").Value = 1064.4445
$ws.Range("K141").Value = 3193.3335
$ws.Range("M141").Value = 1986.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 907.5714
$ws.Range("I50").Value = 1400
$ws.Range("J50").Value = 251
$ws.Range("K50").Value = 1400
$ws.Range("L50").Value = 251
$ws.Range("M50").Value = -686
$ws.Range("N50").Value = -1679

$ws.Range("H61").Value = 125251900
$ws.Range("I61").Value = 143002180
$ws.Range("J61").Value = 1000000
$ws.Range("K61").Value = 143002180
$ws.Range("L61").Value = 1000000
$ws.Range("M61").Value = -143001968
$ws.Range("N61").Value = -1000424

$ws.Range("H63").Value = 3821.3572
$ws.Range("I63").Value = 3853.7693
$ws.Range("J63").Value = 3400
$ws.Range("K63").Value = 3853.7693
$ws.Range("L63").Value = 3400
$ws.Range("M63").Value = -3167.7693
$ws.Range("N63").Value = -4772

$ws.Range("H66").Value = 3821.3572
$ws.Range("I66").Value = 3853.7693
$ws.Range("J66").Value = 3400
$ws.Range("K66").Value = 19268.8465
$ws.Range("L66").Value = 17000
$ws.Range("M66").Value = -15836.8465
$ws.Range("N66").Value = -23864

$ws.Range("H122").Value = 13891568
$ws.Range("I122").Value = 2240
$ws.Range("J122").Value = 22225166
$ws.Range("K122").Value = 6720
$ws.Range("L122").Value = 66675498
$ws.Range("M122").Value = -4270
$ws.Range("N122").Value = -66680398

$ws.Range("H124").Value = 33333.332
$ws.Range("J124").Value = 33333.332
$ws.Range("L124").Value = 33333.332
$ws.Range("N124").Value = -43153.332

$ws.Range("H125").Value = 52741.934
$ws.Range("J125").Value = 52741.934
$ws.Range("L125").Value = 52741.934
$ws.Range("N125").Value = -62581.934

$ws.Range("H132").Value = 8369465
$ws.Range("I132").Value = 9110299
$ws.Range("J132").Value = 220288.6
$ws.Range("K132").Value = 27330897
$ws.Range("L132").Value = 660865.8
$ws.Range("M132").Value = -27328367
$ws.Range("N132").Value = -665925.8

$ws.Range("H136").Value = 125251900
$ws.Range("I136").Value = 143002180
$ws.Range("J136").Value = 1000000
$ws.Range("K136").Value = 429006540
$ws.Range("L136").Value = 3000000
$ws.Range("M136").Value = -429003990
$ws.Range("N136").Value = -3005100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 17456
$ws.Range("J81").Value = 17456
$ws.Range("L81").Value = 17456
$ws.Range("N81").Value = -19578

$ws.Range("H84").Value = 17456
$ws.Range("J84").Value = 17456
$ws.Range("L84").Value = 52368
$ws.Range("N84").Value = -62976

$ws.Range("H105").Value = 38464690
$ws.Range("I105").Value = 83335540
$ws.Range("J105").Value = 3961.5715
$ws.Range("K105").Value = 83335540
$ws.Range("L105").Value = 3961.5715
$ws.Range("M105").Value = -83333793
$ws.Range("N105").Value = -7455.5715

$ws.Range("H135").Value = 52505.363
$ws.Range("J135").Value = 52505.363
$ws.Range("L135").Value = 52505.363
$ws.Range("N135").Value = -62645.363

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 385326.34
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 385326.34
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 385326.34
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -385916.34

$ws.Range("H34").Value = 385326.34
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 385326.34
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 385326.34
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -385730.34

$ws.Range("H105").Value = 1850.4
$ws.Range("I105").Value = 1831.421
$ws.Range("K105").Value = 1831.421
$ws.Range("M105").Value = -84.42100000000005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 242.21875
$ws.Range("I7").Value = 174.25
$ws.Range("J7").Value = 283
$ws.Range("K7").Value = 522.75
$ws.Range("L7").Value = 849
$ws.Range("M7").Value = -410.75
$ws.Range("N7").Value = -1073

$ws.Range("H75").Value = 4958.7144
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4958.7144
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 14876.1432
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -16872.1432

$ws.Range("H78").Value = 4958.7144
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4958.7144
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 44628.4296
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -54612.4296

$ws.Range("H113").Value = 656.8261
$ws.Range("I113").Value = 589.3125
$ws.Range("J113").Value = 811.1429000000001
$ws.Range("K113").Value = 1767.9375
$ws.Range("L113").Value = 2433.4287
$ws.Range("M113").Value = 402.0625
$ws.Range("N113").Value = -6773.4287

$ws.Range("H117").Value = 4167753
$ws.Range("I117").Value = 691.6
$ws.Range("J117").Value = 6061872
$ws.Range("K117").Value = 2074.8
$ws.Range("L117").Value = 18185616
$ws.Range("M117").Value = 1367.2
$ws.Range("N117").Value = -18192500

$ws.Range("H121").Value = 65861370
$ws.Range("J121").Value = 92205760
$ws.Range("L121").Value = 276617280
$ws.Range("N121").Value = -276619900

$ws.Range("H129").Value = 2527089.8
$ws.Range("I129").Value = 1534.3334
$ws.Range("J129").Value = 3474173
$ws.Range("K129").Value = 4603.0002
$ws.Range("L129").Value = 10422519
$ws.Range("M129").Value = 396.9997999999996
$ws.Range("N129").Value = -10432519

$ws.Range("H131").Value = 845.4783
$ws.Range("I131").Value = 435.85715
$ws.Range("J131").Value = 1024.6875
$ws.Range("K131").Value = 1307.57145
$ws.Range("L131").Value = 3074.0625
$ws.Range("M131").Value = 3732.42855
$ws.Range("N131").Value = -13154.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1978.1765
$ws.Range("I122").Value = 1447.3636
$ws.Range("J122").Value = 2951.3333
$ws.Range("K122").Value = 4342.0908
$ws.Range("L122").Value = 8853.999899999999
$ws.Range("M122").Value = -1892.0908
$ws.Range("N122").Value = -13753.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31251788
$ws.Range("I7").Value = 50001840
$ws.Range("J7").Value = 1700
$ws.Range("K7").Value = 50001840
$ws.Range("L7").Value = 1700
$ws.Range("M7").Value = -50001728
$ws.Range("N7").Value = -1924

$ws.Range("H122").Value = 4592.6113
$ws.Range("I122").Value = 6112
$ws.Range("J122").Value = 3832.9167
$ws.Range("K122").Value = 18336
$ws.Range("L122").Value = 11498.7501
$ws.Range("M122").Value = -15886
$ws.Range("N122").Value = -16398.7501

$ws.Range("H126").Value = 31251788
$ws.Range("I126").Value = 50001840
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 150005520
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -150003050
$ws.Range("N126").Value = -10040

$ws.Range("H132").Value = 52497.906
$ws.Range("I132").Value = 3341.0625
$ws.Range("J132").Value = 209799.8
$ws.Range("K132").Value = 10023.1875
$ws.Range("L132").Value = 629399.3999999999
$ws.Range("M132").Value = -7493.1875
$ws.Range("N132").Value = -634459.3999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 37249.355
$ws.Range("I132").Value = 33244.227
$ws.Range("J132").Value = 42215.72
$ws.Range("K132").Value = 99732.681
$ws.Range("L132").Value = 126647.16
$ws.Range("M132").Value = -97202.681
$ws.Range("N132").Value = -131707.16

Write-Host "Hades_Profits refresh applied across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR"
